# Apply weekly price update for Fruta / hortaliza (Chirimoya) sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44839
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 160
$ws.Range("N2").Value = 26000
$ws.Range("O2").Value = 27000
$ws.Range("P2").Value = 26500
$ws.Range("R2").Value = "Región de Coquimbo"
$ws.Range("S2").Value = 2208

# Row 3
$ws.Range("D3").Value = 44482
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 160
$ws.Range("N3").Value = 25000
$ws.Range("O3").Value = 26000
$ws.Range("P3").Value = 25500
$ws.Range("S3").Value = 2125

# Row 4
$ws.Range("D4").Value = 44545
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 23000
$ws.Range("O4").Value = 24000
$ws.Range("P4").Value = 23500
$ws.Range("Q4").Value = "`$/bandeja 12 kilos"
$ws.Range("S4").Value = 1958

# Row 7
$ws.Range("D7").Value = 44783
$ws.Range("L7").Value = "Tercera"
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 27000
$ws.Range("O7").Value = 28000
$ws.Range("P7").Value = 27500
$ws.Range("S7").Value = 2292

# Row 9
$ws.Range("D9").Value = 44160
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 19000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 19500
$ws.Range("Q9").Value = "`$/caja 13 kilos"
$ws.Range("S9").Value = 1500
$ws.Range("T9").Value = 13

# Row 10
$ws.Range("D10").Value = 44167
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 18000
$ws.Range("O10").Value = 19000
$ws.Range("P10").Value = 18500
$ws.Range("Q10").Value = "`$/caja 13 kilos"
$ws.Range("S10").Value = 1423
$ws.Range("T10").Value = 13

# Row 11
$ws.Range("D11").Value = 44776
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 160
$ws.Range("Q11").Value = "`$/caja 10 kilos"
$ws.Range("S11").Value = 2950
$ws.Range("T11").Value = 10

# Row 12
$ws.Range("D12").Value = 44489
$ws.Range("N12").Value = 24000
$ws.Range("O12").Value = 25000
$ws.Range("P12").Value = 24500
$ws.Range("Q12").Value = "`$/caja 12 kilos"
$ws.Range("S12").Value = 2042
$ws.Range("T12").Value = 12

# Row 13
$ws.Range("D13").Value = 44468
$ws.Range("N13").Value = 29000
$ws.Range("O13").Value = 30000
$ws.Range("P13").Value = 29500
$ws.Range("Q13").Value = "`$/bandeja 10 kilos"
$ws.Range("S13").Value = 2950
$ws.Range("T13").Value = 10

# Row 14
$ws.Range("D14").Value = 44475
$ws.Range("L14").Value = "Especial"
$ws.Range("N14").Value = 32000
$ws.Range("O14").Value = 33000
$ws.Range("P14").Value = 32500
$ws.Range("S14").Value = 2708

# Row 16
$ws.Range("D16").Value = 44881
$ws.Range("L16").Value = "Primera"
$ws.Range("N16").Value = 22000
$ws.Range("O16").Value = 23000
$ws.Range("P16").Value = 22500
$ws.Range("Q16").Value = "`$/caja 12 kilos"
$ws.Range("S16").Value = 1875
$ws.Range("T16").Value = 12

# Row 17
$ws.Range("D17").Value = 44860
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 23000
$ws.Range("O17").Value = 24000
$ws.Range("P17").Value = 23500
$ws.Range("R17").Value = "Provincia de Limarí"
$ws.Range("S17").Value = 1958

# Row 18
$ws.Range("D18").Value = 44811
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 29000
$ws.Range("O18").Value = 30000
$ws.Range("P18").Value = 29500
$ws.Range("Q18").Value = "`$/caja 12 kilos"
$ws.Range("S18").Value = 2458
$ws.Range("T18").Value = 12

# Row 20
$ws.Range("D20").Value = 44496
$ws.Range("N20").Value = 23000
$ws.Range("O20").Value = 24000
$ws.Range("P20").Value = 23500
$ws.Range("S20").Value = 1958
